# Edit: "Updates to Measurment Chain + Current Calc"
# Adds two new worksheets ("Battery Estimate" and "Load Cell Current") between
# the existing "INA + VGA" and "nur INA" sheets, with their data/formulas,
# plus small cosmetic updates (column widths, wrap-text header, selections).

$wb = $excel.ActiveWorkbook

$inaVga = $wb.Worksheets.Item("INA + VGA")
$nurIna = $wb.Worksheets.Item("nur INA")

# xlVAlignCenter
$xlVAlignCenter = -4108

# --- Create sheets in an order that reproduces the target sheetId numbering:
# "Load Cell Current" is created first (sheetId=3), then "Battery Estimate" is
# inserted right after "INA + VGA" (sheetId=4), which pushes "Load Cell
# Current" to the 3rd tab position. Final tab order:
#   INA + VGA, Battery Estimate, Load Cell Current, nur INA
#
# NOTE: object references returned by Worksheets.Add() go stale once another
# sheet is added afterwards, so every sheet used below is re-fetched by name
# via Worksheets.Item(...) right before it is used.
$tmp = $wb.Worksheets.Add($null, $inaVga)
$tmp.Name = "Load Cell Current"

$tmp2 = $wb.Worksheets.Add($null, $inaVga)
$tmp2.Name = "Battery Estimate"

$lcc = $wb.Worksheets.Item("Load Cell Current")
$be = $wb.Worksheets.Item("Battery Estimate")

# =====================================================================
# Text labels first, entered in the same order the original author typed
# them in (keeps the shared-strings table in the same order as the source).
# =====================================================================

# -- Load Cell Current headers --
$lcc.Range("C1").Value = "I @5V [mA]"
$lcc.Range("D1").Value = "I @10V [mA]"
$lcc.Range("B1").Value = "Bridge Resistance [Ohm]"
$lcc.Range("A1").Value = "Load Cell [kgf]"

# -- Battery Estimate headers + labels --
$be.Range("A1").Value = "Component"
$be.Range("B1").Value = "I_typ [mA]"
$be.Range("C1").Value = "I_max [mA]"
$be.Range("G4").Value = "Capacity"
$be.Range("A2").Value = "Load Cell"
$be.Range("A3").Value = "LDO"
$be.Range("A4").Value = "INA"
$be.Range("A5").Value = "OPAMP"
$be.Range("A6").Value = "Digi Poti"
$be.Range("G6").Value = "Estimate Typ [h]"
$be.Range("G5").Value = "Estimate Min [h]"
$be.Range("A7").Value = "DS ADC"
$be.Range("G2").Value = "Total I_typ"
$be.Range("G3").Value = "Total I_max"
$be.Range("A8").Value = "VREF"
$be.Range("J2").Value = "Total I_typ"
$be.Range("J3").Value = "Total I_max"

# =====================================================================
# Battery Estimate - numeric data + formulas
# =====================================================================

$be.Range("B2").Value = 14.28
$be.Range("C2").Value = 28.57

$be.Range("B3").Value = 5
$be.Range("C3").Value = 6

$be.Range("B4").Value = 2
$be.Range("C4").Value = 2.3

$be.Range("B5").Value = 0.45
$be.Range("C5").Value = 0.51

$be.Range("B6").Value = 0.001
$be.Range("C6").Value = 0.002

$be.Range("B7").Value = 0.145
$be.Range("C7").Value = 0.16

$be.Range("B8").Value = 0.38
$be.Range("C8").Value = 0.58

$be.Range("H2").Formula = "=SUM(B2:B100)"
$be.Range("K2").Formula = "=SUM(B3:B100)"

$be.Range("H3").Formula = "=SUM(C2:C100)"
$be.Range("K3").Formula = "=SUM(C3:C100)"

$be.Range("H4").Value = 5000

$be.Range("H5").Formula = "=H4/H3"
$be.Range("H6").Formula = "=H4/H2"

# Column widths (COM ColumnWidth ~= stored xlsx width - 5/6)
$be.Columns.Item(1).ColumnWidth = 13 - 5/6
$be.Columns.Item(2).ColumnWidth = 10.88671875 - 5/6
$be.Columns.Item(3).ColumnWidth = 11.77734375 - 5/6
$be.Columns.Item(7).ColumnWidth = 13.6640625 - 5/6

# =====================================================================
# Load Cell Current - numeric data + formulas
# =====================================================================

# Column A uses the same "vertical-center + wrap text" look as the Load Cell
# column on the other two sheets.
$lcc.Range("A1:A7").WrapText = $true
$lcc.Range("A1:A7").VerticalAlignment = $xlVAlignCenter

# Header for the resistance column additionally wraps (long header, narrow column).
$lcc.Range("B1").WrapText = $true

$lcc.Range("A2").Value = 0.15
$lcc.Range("B2").Value = 350

$lcc.Range("A3").Value = 2
$lcc.Range("B3").Value = 350

$lcc.Range("A4").Value = 20
$lcc.Range("B4").Value = 175

$lcc.Range("A5").Value = 45
$lcc.Range("B5").Value = 350

$lcc.Range("A6").Value = 250
$lcc.Range("B6").Value = 175

$lcc.Range("A7").Value = 3000
$lcc.Range("B7").Value = 175

# Formulas (C2:C7 and D2:D7 filled together so they come out as shared formulas)
$lcc.Range("C2:C7").Formula = '=5/$B2 * 1000'
$lcc.Range("D2:D7").Formula = '=10/$B2 * 1000'

# Column widths
$lcc.Columns.Item(2).ColumnWidth = 15.21875 - 5/6
$lcc.Columns.Item(3).ColumnWidth = 12 - 5/6
$lcc.Columns.Item(4).ColumnWidth = 11.109375 - 5/6

# =====================================================================
# Sheet-view / selection tweaks
# =====================================================================

# "INA + VGA": selection becomes A1:A7, no longer the tab-selected sheet
$inaVga.Range("A1:A7").Select()

# "nur INA": keep its own previous selection (J4), no longer tab-selected
$nurIna.Range("J4").Select()

# "Load Cell Current": whole header row selected
$lcc.Range("A1:XFD1").Select()

# "Battery Estimate" ends up the active / tab-selected sheet, cursor at H11
$be.Activate()
$be.Range("H11").Select()
